$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 230, shifting the existing block (230-237) down to (232-239)
$ws.Rows.Item(230).Resize(2).Insert()

# New row 230: Feria Lagunitas de Puerto Montt - Palta Hass "Primera", week of 2021-11-09
$ws.Cells.Item(230, 1).Value = 4
$ws.Cells.Item(230, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(230, 3).Value = "Los Lagos"
$ws.Cells.Item(230, 4).Value = 44509
$ws.Cells.Item(230, 5).Value = 10
$ws.Cells.Item(230, 6).Value = "Fruta"
$ws.Cells.Item(230, 7).Value = 100106
$ws.Cells.Item(230, 8).Value = "Oleaginosos"
$ws.Cells.Item(230, 9).Value = 100106002
$ws.Cells.Item(230, 10).Value = "Palta"
$ws.Cells.Item(230, 11).Value = "Hass"
$ws.Cells.Item(230, 12).Value = "Primera"
$ws.Cells.Item(230, 13).Value = 400
$ws.Cells.Item(230, 14).Value = 4000
$ws.Cells.Item(230, 15).Value = 4200
$ws.Cells.Item(230, 16).Value = 4100
$ws.Cells.Item(230, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(230, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(230, 19).Value = 4100
$ws.Cells.Item(230, 20).Value = 1

# New row 231: Feria Lagunitas de Puerto Montt - Palta Hass "Segunda", week of 2021-11-09
$ws.Cells.Item(231, 1).Value = 4
$ws.Cells.Item(231, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(231, 3).Value = "Los Lagos"
$ws.Cells.Item(231, 4).Value = 44509
$ws.Cells.Item(231, 5).Value = 10
$ws.Cells.Item(231, 6).Value = "Fruta"
$ws.Cells.Item(231, 7).Value = 100106
$ws.Cells.Item(231, 8).Value = "Oleaginosos"
$ws.Cells.Item(231, 9).Value = 100106002
$ws.Cells.Item(231, 10).Value = "Palta"
$ws.Cells.Item(231, 11).Value = "Hass"
$ws.Cells.Item(231, 12).Value = "Segunda"
$ws.Cells.Item(231, 13).Value = 200
$ws.Cells.Item(231, 14).Value = 3600
$ws.Cells.Item(231, 15).Value = 3600
$ws.Cells.Item(231, 16).Value = 3600
$ws.Cells.Item(231, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(231, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(231, 19).Value = 3600
$ws.Cells.Item(231, 20).Value = 1

# Format the date cells in the newly inserted rows the same way as the rest of column D
$ws.Range("D230:D231").NumberFormat = $ws.Range("D232").NumberFormat
